$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Matteo Mazzola"
$ws.Range("B6").Value = "MATTEO PILATI | Pinguini Trentini"
$ws.Range("C6").Value = "Leonardo  Parisi  | MediaserT"
$ws.Range("D6").Value = "Andrea Conzatti | FC Savignano"
$ws.Range("E6").Value = "Andrea Riolfatti | La Contea FC"
$ws.Range("F6").Value = "Francesco Cristoforetti | Vigili del Fusto"
